# Apply updated crypto price/volume data per GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.492.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.80%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.490.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.489.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.33%  "

# Row 10
$ws.Range("E10").Value = "  -2.43%  "

# Row 11
$ws.Range("E11").Value = "  -0.42%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.949.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.406.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.523.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "

# Row 22
$ws.Range("E22").Value = "  -1.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.11%  "

# Row 24
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.29%  "

# Row 26
$ws.Range("E26").Value = "  -2.78%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.621.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0871"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.07%  "

# Row 32
$ws.Range("E32").Value = "  -6.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "437.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("E35").Value = "  -2.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.91%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.40%  "

# Row 40
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.45%  "

# Row 43
$ws.Range("E43").Value = "  -2.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.40%  "

# Row 45
$ws.Range("E45").Value = "  -7.55%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "138.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.511"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.73%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0722"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.37%  "

# Row 50
$ws.Range("E50").Value = "  -1.16%  "

# Row 51
$ws.Range("B51").Value = "POPCAT"
$ws.Range("C51").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +20.53%  "
